$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two product rows: row 2 becomes the BS-2DQRUSBWS entry with 8 units,
# row 3 becomes the ROPA001 entry with 9 units.
$ws.Range("A2").Value = "BS-2DQRUSBWS"
$ws.Range("B2").Value = 8

$ws.Range("A3").Value = "ROPA001"
$ws.Range("B3").Value = 9
